$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text cells (id, phone number, name, date-as-text) to keep their
# literal text representation instead of being auto-coerced by Excel into
# a number (phone) or a date serial (date string). We temporarily apply a
# text number format, assign the value, then restore the default "Normal"
# style so the cells end up with no explicit style, same as the source.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "75c44810a32a3d6447df"
$ws.Range("B2").Value = "+74267426016"
$ws.Range("C2").Value = "Automation User 10"
$ws.Range("E2").Value = "2026-02-19"

$ws.Range("A2").Style = "Normal"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Style = "Normal"
$ws.Range("E2").Style = "Normal"

# Payment amount becomes a genuine numeric value (no longer text).
$ws.Range("D2").Value = 100
